# ----------------------------------------------------------------------------
# Applies the "semi final" edit to the checkout / active_tabs workbook:
#  - checkout (sheet1) is emptied out (its header row content is gone)
#  - active_tabs (sheet2) is rebuilt with the new transactional order data,
#    a reworded/re-styled header, and a tweaked column width
#  - checkout becomes the active / selected sheet instead of active_tabs
# ----------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "checkout"
$ws2 = $wb.Worksheets.Item(2)   # "active_tabs"

# -----------------------------------------------------------------
# 1) active_tabs: create an orphaned center-aligned style (matches
#    the extra unused cellXf that shows up at index 2 in the target
#    styles table) by touching a throw-away cell and then discarding
#    the row it lived in.
# -----------------------------------------------------------------
$scratch = $ws2.Range("A50")
$scratch.HorizontalAlignment = -4108   # xlCenter
$ws2.Rows.Item(50).Delete() | Out-Null

# -----------------------------------------------------------------
# 2) active_tabs: wipe the old order rows, but keep the header row's
#    existing (bold + centered) formatting in place
# -----------------------------------------------------------------
$ws2.Range("A2:E9").Clear() | Out-Null
$ws2.Range("C1").Clear() | Out-Null   # drop the stray empty numeric header cell

# Header row text
$ws2.Range("A1").Value2 = "table Number"
$ws2.Range("B1").Value2 = "Orders"

# Give B1 the bold + centered + custom date-time number format style
# (this creates / reuses the numFmt 164 "yyyy-mm-dd h:mm:ss" entry and
# the matching cellXf, landing on index 3 right after the orphan above)
$ws2.Range("B1").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"

# New order rows. Column A values look numeric ("03", "10", ...) but
# must stay textual (leading zeros preserved), so they are entered as
# formulas producing text and then frozen into plain values.
$tableNumbers = @("03", "03", "10", "02", "20", "05")
for ($i = 0; $i -lt $tableNumbers.Length; $i++) {
    $row = 2 + $i
    $ws2.Range("A$row").Formula = '="' + $tableNumbers[$i] + '"'
}
$numRng = $ws2.Range("A2:A7")
$numRng.Copy() | Out-Null
$numRng.PasteSpecial(-4163) | Out-Null   # xlPasteValues

$ws2.Range("B2").Value2 = "09/07/2023 20:02:15"
$ws2.Range("C2").Value2 = "ice tea"

$ws2.Range("B3").Value2 = "09/07/2023 20:02:23"
$ws2.Range("C3").Value2 = "ice cream"

$ws2.Range("B4").Value2 = "09/07/2023 20:03:58"
$ws2.Range("C4").Value2 = "ice 3 in 1"
$ws2.Range("D4").Value2 = "ice coffee"

$ws2.Range("B5").Value2 = "09/07/2023 20:24:16"
$ws2.Range("C5").Value2 = "tea"

$ws2.Range("B6").Value2 = "09/07/2023 20:24:22"
$ws2.Range("C6").Value2 = "coffee"

$ws2.Range("B7").Value2 = "09/07/2023 20:24:33"
$ws2.Range("C7").Value2 = "3 in 1"
$ws2.Range("D7").Value2 = "ice 3 in 1"

# Column B grows a little to fit the new date/time strings
$ws2.Columns.Item(2).ColumnWidth = 18.109375

# -----------------------------------------------------------------
# 3) checkout: clear out its old header row entirely
# -----------------------------------------------------------------
$ws1.Range("A1:B1").Clear() | Out-Null

# -----------------------------------------------------------------
# 4) Page setup: both sheets end up with an explicit portrait setup
# -----------------------------------------------------------------
$ws1.PageSetup.Orientation = 1   # xlPortrait
$ws2.PageSetup.Orientation = 1   # xlPortrait

# -----------------------------------------------------------------
# 5) View state: checkout becomes the active / selected sheet, with
#    its own selection, while active_tabs keeps a plain selection.
# -----------------------------------------------------------------
$ws2.Range("D4").Select() | Out-Null
$ws2.Range("A2:D4").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("A3").Select() | Out-Null
$ws1.Range("A1:J3").Select() | Out-Null

Write-Host "edit complete"
